$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows after the existing row 217, shifting rows 218:229
# down to 220:231 (formatting/styles carried down with the shift).
$ws.Rows("218:219").Insert()

# New row 218: weekly Albahaca "Primera" data point for Provincia de Chacabuco
$ws.Cells.Item(218, 1).Value = 9
$ws.Cells.Item(218, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(218, 3).Value = "Metropolitana"
$ws.Cells.Item(218, 4).Value = 44516
$ws.Cells.Item(218, 5).Value = 13
$ws.Cells.Item(218, 6).Value = 100112052
$ws.Cells.Item(218, 7).Value = "Albahaca"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 52
$ws.Cells.Item(218, 11).Value = 6000
$ws.Cells.Item(218, 12).Value = 7000
$ws.Cells.Item(218, 13).Value = 6500
$ws.Cells.Item(218, 14).Value = "$/docena de matas"
$ws.Cells.Item(218, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(218, 16).Value = 1083
$ws.Cells.Item(218, 17).Value = 6
$ws.Cells.Item(218, 18).Value = "Hortaliza"

# New row 219: weekly Albahaca "Primera" data point for Provincia del Elquí
$ws.Cells.Item(219, 1).Value = 9
$ws.Cells.Item(219, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(219, 3).Value = "Metropolitana"
$ws.Cells.Item(219, 4).Value = 44516
$ws.Cells.Item(219, 5).Value = 13
$ws.Cells.Item(219, 6).Value = 100112052
$ws.Cells.Item(219, 7).Value = "Albahaca"
$ws.Cells.Item(219, 8).Value = "Sin especificar"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 25
$ws.Cells.Item(219, 11).Value = 7000
$ws.Cells.Item(219, 12).Value = 8000
$ws.Cells.Item(219, 13).Value = 7520
$ws.Cells.Item(219, 14).Value = "$/docena de matas"
$ws.Cells.Item(219, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(219, 16).Value = 1253
$ws.Cells.Item(219, 17).Value = 6
$ws.Cells.Item(219, 18).Value = "Hortaliza"
